$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "helper table" that lived in columns J:U, rows 2-8.
# It was a scratch area used to plan shifts (names / letter codes / counts)
# that is being replaced by the new in-GUI "create table" feature.
$ws.Range("J2:U8").Clear()

# J9 ("תגבור") keeps its text but loses the highlighted fill style it used
# to carry; instead J2 becomes an empty cell that carries that highlight
# style, ready for the new table feature.
$ws.Range("J9").Style = "Normal"

$ws.Range("J2").Value = $null
$ws.Range("J2").Interior.Color = 10025880
